$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 189 (shifts old rows 189-192 down to 191-194)
$ws.Rows.Item(189).Insert()
$ws.Rows.Item(189).Insert()

# New row 189: Caqui / Mankaki / Primera, dated 2023-05-31 (serial 45077)
$ws.Range("A189").Value = 6
$ws.Range("B189").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C189").Value = "Metropolitana"
$ws.Range("D189").Value = 45077
$ws.Range("D189").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E189").Value = 13
$ws.Range("F189").Value = "Fruta"
$ws.Range("G189").Value = 100107
$ws.Range("H189").Value = "Otros"
$ws.Range("I189").Value = 100107001
$ws.Range("J189").Value = "Caqui"
$ws.Range("K189").Value = "Mankaki"
$ws.Range("L189").Value = "Primera"
$ws.Range("M189").Value = 15
$ws.Range("N189").Value = 260000
$ws.Range("O189").Value = 260000
$ws.Range("P189").Value = 260000
$ws.Range("Q189").Value = "$/bins (450 kilos)"
$ws.Range("R189").Value = "Región de O'Higgins"
$ws.Range("S189").Value = 578
$ws.Range("T189").Value = 450

# New row 190: Caqui / Mankaki / Segunda, dated 2023-05-31 (serial 45077)
$ws.Range("A190").Value = 6
$ws.Range("B190").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C190").Value = "Metropolitana"
$ws.Range("D190").Value = 45077
$ws.Range("D190").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E190").Value = 13
$ws.Range("F190").Value = "Fruta"
$ws.Range("G190").Value = 100107
$ws.Range("H190").Value = "Otros"
$ws.Range("I190").Value = 100107001
$ws.Range("J190").Value = "Caqui"
$ws.Range("K190").Value = "Mankaki"
$ws.Range("L190").Value = "Segunda"
$ws.Range("M190").Value = 15
$ws.Range("N190").Value = 230000
$ws.Range("O190").Value = 230000
$ws.Range("P190").Value = 230000
$ws.Range("Q190").Value = "$/bins (450 kilos)"
$ws.Range("R190").Value = "Región de O'Higgins"
$ws.Range("S190").Value = 511
$ws.Range("T190").Value = 450

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
